# Populate the "Courses" worksheet with the course-info header row and one
# sample data row, matching the courseTemp.xlsx upload format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) — all text
$ws.Range("A1").Value = "department"
$ws.Range("B1").Value = "number"
$ws.Range("C1").Value = "univNumber"
$ws.Range("D1").Value = "name"
$ws.Range("E1").Value = "category"
$ws.Range("F1").Value = "topic"
$ws.Range("G1").Value = "hours"
$ws.Range("H1").Value = "section"
$ws.Range("I1").Value = "faculty"
$ws.Range("J1").Value = "semester"

# Data row (row 2)
$ws.Range("A2").Value = "COMP"

# "number" / "section" look numeric but are stored as text in the source
# data, so force text entry (leading apostrophe) and strip the resulting
# quote-prefix style back to Normal so no extra formatting is introduced.
$ws.Range("B2").Value = "'101"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = 11
$ws.Range("D2").Value = "Fluency in Information Technology"
$ws.Range("E2").Value = "Theory"
$ws.Range("F2").Value = "N/A"
$ws.Range("G2").Value = 3

$ws.Range("H2").Value = "'1"
$ws.Range("H2").Style = "Normal"

$ws.Range("I2").Value = "Ahalt, Stanley"
$ws.Range("J2").Value = "FA 2018"
